$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.389.03'
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").Value = '''1.848.40'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("D4").Value = '''0.9993'
$ws.Range("D5").Value = '''241.18'
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").Value = '''0.6254'
$ws.Range("E6").Value = '  -3.97%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '''0.07617'
$ws.Range("E8").Value = '  +1.55%  '
$ws.Range("D9").Value = '''0.2968'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").Value = '''24.40'
$ws.Range("E10").Value = '  -0.34%  '
$ws.Range("D11").Value = '''2.061.75'
$ws.Range("E11").Value = '  +11.38%  '
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("D13").Value = '''4.990'
$ws.Range("E13").Value = '  -1.00%  '
$ws.Range("D14").Value = '''0.6877'
$ws.Range("E14").Value = '  +0.13%  '
$ws.Range("D15").Value = '''82.98'
$ws.Range("E15").Value = '  -0.74%  '
$ws.Range("D16").Value = '''0.000009940'
$ws.Range("E16").Value = '  +4.09%  '
$ws.Range("D17").Value = '''2.261.71'
$ws.Range("E17").Value = '  +7.46%  '
$ws.Range("D18").Value = '''6.156'
$ws.Range("E18").Value = '  +0.21%  '
$ws.Range("D19").Value = '''29.679.63'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").Value = '''230.83'
$ws.Range("E20").Value = '  -2.72%  '
$ws.Range("D21").Value = '''12.53'
$ws.Range("E21").Value = '  -0.76%  '
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '''7.632'
$ws.Range("E23").Value = '  -1.01%  '
$ws.Range("D24").Value = '''1.001'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").Value = '''154.66'
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D27").Value = '''8.466'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  -0.92%  '
$ws.Range("D29").Value = '''1.475'
$ws.Range("E29").Value = '  -1.12%  '
$ws.Range("D30").Value = '''0.05812'
$ws.Range("E30").Value = '  -4.35%  '
$ws.Range("D31").Value = '''1.255'
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '''4.124'
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("D33").Value = '''4.017'
$ws.Range("E33").Value = '  -1.51%  '
$ws.Range("D34").Value = '''1.872'
$ws.Range("E34").Value = '  +0.06%  '
$ws.Range("D35").Value = '''1.160'
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("D36").Value = '''0.7173'
$ws.Range("E36").Value = '  -1.26%  '
$ws.Range("D37").Value = '''2.597'
$ws.Range("E37").Value = '  -0.23%  '
$ws.Range("D38").Value = '''1.251.61'
$ws.Range("E38").Value = '  +4.11%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '''0.01802'
$ws.Range("E40").Value = '  +0.97%  '
$ws.Range("B41").Value = 'RocketPoolETH'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D41").Value = '''2.201.41'
$ws.Range("E41").Value = '  +9.41%  '
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = '''0.9073'
$ws.Range("E42").Value = '  -0.17%  '
$ws.Range("D43").Value = '''6.080'
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("D44").Value = '''0.9998'
$ws.Range("D45").Value = '''102.01'
$ws.Range("E45").Value = '  +0.45%  '
$ws.Range("D46").Value = '''67.49'
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("D47").Value = '''7.318'
$ws.Range("E47").Value = '  -0.61%  '
$ws.Range("E48").Value = '  -2.78%  '
$ws.Range("D49").Value = '''9.134'
$ws.Range("E49").Value = '  +0.08%  '
$ws.Range("D50").Value = '''0.4025'
$ws.Range("E50").Value = '  -0.74%  '
$ws.Range("E51").Value = '  +2.11%  '
